# Preparation for publication 0.2.0
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1. Bump the published version number
$ws.Range("B3").Value = "0.2.0"

# 2. Refresh the publication date/time
$ws.Range("B8").Value = "2023-10-19T17:05:12+00:00"

# 3. Insert a new "Jurisdiction" row right after "Contact" (row 10), which
#    pushes Description/Purpose/Copyright/FHIR Version/etc. down by one row.
$ws.Rows("11:11").Insert()
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"

# Match the formatting used by the rest of the body rows (the inserted row
# otherwise comes through with the default/no style).
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false
